$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 18 and row 19 (columns B:AD), keep column A (index) unchanged
$ws.Range("B18").Value2 = 7003585
$ws.Range("B19").Value2 = 7004591
$ws.Range("C18").Value2 = 'Qatar Stars League'
$ws.Range("C19").Value2 = 'Qatar Stars League'
$ws.Range("D18").Value2 = 45171.59375
$ws.Range("D19").Value2 = 45171.59375
$ws.Range("E18").Value2 = 'Al Sadd'
$ws.Range("E19").Value2 = 'AlShamal SC'
$ws.Range("F18").Value2 = 'AlWakrah SC'
$ws.Range("F19").Value2 = 'AlRayyan SC'
$ws.Range("G18").Value2 = 0
$ws.Range("G19").Value2 = 3
$ws.Range("H18").Value2 = 0
$ws.Range("H19").Value2 = 4
$ws.Range("I18").Value2 = 0
$ws.Range("I19").Value2 = 1
$ws.Range("J18").Value2 = 0
$ws.Range("J19").Value2 = 2
$ws.Range("K18").Value2 = 'D'
$ws.Range("K19").Value2 = 'A'
$ws.Range("L18").Value2 = 1.615
$ws.Range("L19").Value2 = 4.5
$ws.Range("M18").Value2 = 4
$ws.Range("M19").Value2 = 4.2
$ws.Range("N18").Value2 = 4.333
$ws.Range("N19").Value2 = 1.55
$ws.Range("O18").Value2 = 1.533
$ws.Range("O19").Value2 = 3.3
$ws.Range("P18").Value2 = 4.2
$ws.Range("P19").Value2 = 3.8
$ws.Range("Q18").Value2 = 5
$ws.Range("Q19").Value2 = 1.85
$ws.Range("R18").Value2 = -1
$ws.Range("R19").Value2 = 0.5
$ws.Range("S18").Value2 = 1.8
$ws.Range("S19").Value2 = 1.85
$ws.Range("T18").Value2 = 2
$ws.Range("T19").Value2 = 1.95
$ws.Range("U18").Value2 = 3.5
$ws.Range("U19").Value2 = 2.75
$ws.Range("V18").Value2 = 1.925
$ws.Range("V19").Value2 = 1.85
$ws.Range("W18").Value2 = 1.875
$ws.Range("W19").Value2 = 1.95
$ws.Range("X18").Value2 = -1
$ws.Range("X19").Value2 = -1
$ws.Range("Y18").Value2 = 3.2
$ws.Range("Y19").Value2 = -1
$ws.Range("Z18").Value2 = -1
$ws.Range("Z19").Value2 = 0.8500000000000001
$ws.Range("AA18").Value2 = -1
$ws.Range("AA19").Value2 = -1
$ws.Range("AB18").Value2 = 1
$ws.Range("AB19").Value2 = 0.95
$ws.Range("AC18").Value2 = -1
$ws.Range("AC19").Value2 = 0.8500000000000001
$ws.Range("AD18").Value2 = 0.875
$ws.Range("AD19").Value2 = -1

# Swap row 21 and row 22 (columns B:AD), keep column A (index) unchanged
$ws.Range("B21").Value2 = 7003586
$ws.Range("B22").Value2 = 7003478
$ws.Range("C21").Value2 = 'Qatar Stars League'
$ws.Range("C22").Value2 = 'Qatar Stars League'
$ws.Range("D21").Value2 = 45192.58333333334
$ws.Range("D22").Value2 = 45192.58333333334
$ws.Range("E21").Value2 = 'Qatar SC Doha'
$ws.Range("E22").Value2 = 'Umm Salal'
$ws.Range("F21").Value2 = 'Al Sadd'
$ws.Range("F22").Value2 = 'AlAhli Doha'
$ws.Range("G21").Value2 = 1
$ws.Range("G22").Value2 = 3
$ws.Range("H21").Value2 = 3
$ws.Range("H22").Value2 = 1
$ws.Range("I21").Value2 = 1
$ws.Range("I22").Value2 = 2
$ws.Range("J21").Value2 = 1
$ws.Range("J22").Value2 = 0
$ws.Range("K21").Value2 = 'A'
$ws.Range("K22").Value2 = 'H'
$ws.Range("L21").Value2 = 5
$ws.Range("L22").Value2 = 2
$ws.Range("M21").Value2 = 4.333
$ws.Range("M22").Value2 = 3.6
$ws.Range("N21").Value2 = 1.55
$ws.Range("N22").Value2 = 3.25
$ws.Range("O21").Value2 = 5
$ws.Range("O22").Value2 = 1.95
$ws.Range("P21").Value2 = 4.2
$ws.Range("P22").Value2 = 3.6
$ws.Range("Q21").Value2 = 1.571
$ws.Range("Q22").Value2 = 3.4
$ws.Range("R21").Value2 = 1
$ws.Range("R22").Value2 = -0.5
$ws.Range("S21").Value2 = 1.85
$ws.Range("S22").Value2 = 1.975
$ws.Range("T21").Value2 = 1.95
$ws.Range("T22").Value2 = 1.825
$ws.Range("U21").Value2 = 3
$ws.Range("U22").Value2 = 3
$ws.Range("V21").Value2 = 1.85
$ws.Range("V22").Value2 = 1.8
$ws.Range("W21").Value2 = 1.95
$ws.Range("W22").Value2 = 2
$ws.Range("X21").Value2 = -1
$ws.Range("X22").Value2 = 0.95
$ws.Range("Y21").Value2 = -1
$ws.Range("Y22").Value2 = -1
$ws.Range("Z21").Value2 = 0.571
$ws.Range("Z22").Value2 = -1
$ws.Range("AA21").Value2 = -1
$ws.Range("AA22").Value2 = 0.9750000000000001
$ws.Range("AB21").Value2 = 0.95
$ws.Range("AB22").Value2 = -1
$ws.Range("AC21").Value2 = 0.8500000000000001
$ws.Range("AC22").Value2 = 0.8
$ws.Range("AD21").Value2 = -1
$ws.Range("AD22").Value2 = -1

# Swap row 24 and row 25 (columns B:AD), keep column A (index) unchanged
$ws.Range("B24").Value2 = 7004592
$ws.Range("B25").Value2 = 7004593
$ws.Range("C24").Value2 = 'Qatar Stars League'
$ws.Range("C25").Value2 = 'Qatar Stars League'
$ws.Range("D24").Value2 = 45193.58333333334
$ws.Range("D25").Value2 = 45193.58333333334
$ws.Range("E24").Value2 = 'AlMuaidar'
$ws.Range("E25").Value2 = 'Al Markhiya'
$ws.Range("F24").Value2 = 'AlShamal SC'
$ws.Range("F25").Value2 = 'AlWakrah SC'
$ws.Range("G24").Value2 = 2
$ws.Range("G25").Value2 = 0
$ws.Range("H24").Value2 = 2
$ws.Range("H25").Value2 = 3
$ws.Range("I24").Value2 = 0
$ws.Range("I25").Value2 = 0
$ws.Range("J24").Value2 = 0
$ws.Range("J25").Value2 = 1
$ws.Range("K24").Value2 = 'D'
$ws.Range("K25").Value2 = 'A'
$ws.Range("L24").Value2 = 2.4
$ws.Range("L25").Value2 = 4
$ws.Range("M24").Value2 = 3.6
$ws.Range("M25").Value2 = 4
$ws.Range("N24").Value2 = 2.4
$ws.Range("N25").Value2 = 1.727
$ws.Range("O24").Value2 = 2.4
$ws.Range("O25").Value2 = 4.75
$ws.Range("P24").Value2 = 3.6
$ws.Range("P25").Value2 = 4.333
$ws.Range("Q24").Value2 = 2.45
$ws.Range("Q25").Value2 = 1.571
$ws.Range("R24").Value2 = 0
$ws.Range("R25").Value2 = 1
$ws.Range("S24").Value2 = 1.925
$ws.Range("S25").Value2 = 1.85
$ws.Range("T24").Value2 = 1.875
$ws.Range("T25").Value2 = 1.95
$ws.Range("U24").Value2 = 2.75
$ws.Range("U25").Value2 = 3
$ws.Range("V24").Value2 = 1.75
$ws.Range("V25").Value2 = 1.825
$ws.Range("W24").Value2 = 1.95
$ws.Range("W25").Value2 = 1.975
$ws.Range("X24").Value2 = -1
$ws.Range("X25").Value2 = -1
$ws.Range("Y24").Value2 = 2.6
$ws.Range("Y25").Value2 = -1
$ws.Range("Z24").Value2 = -1
$ws.Range("Z25").Value2 = 0.571
$ws.Range("AA24").Value2 = 0
$ws.Range("AA25").Value2 = -1
$ws.Range("AB24").Value2 = 0
$ws.Range("AB25").Value2 = 0.95
$ws.Range("AC24").Value2 = 0.75
$ws.Range("AC25").Value2 = 0
$ws.Range("AD24").Value2 = -1
$ws.Range("AD25").Value2 = 0

# Swap row 68 and row 69 (columns B:AD), keep column A (index) unchanged
$ws.Range("B68").Value2 = 7609336
$ws.Range("B69").Value2 = 7004626
$ws.Range("C68").Value2 = 'Qatar Stars League'
$ws.Range("C69").Value2 = 'Qatar Stars League'
$ws.Range("D68").Value2 = 45280.47916666666
$ws.Range("D69").Value2 = 45280.47916666666
$ws.Range("E68").Value2 = 'AlMuaidar'
$ws.Range("E69").Value2 = 'Al Gharafa'
$ws.Range("F68").Value2 = 'AlWakrah SC'
$ws.Range("F69").Value2 = 'Qatar SC Doha'
$ws.Range("G68").Value2 = 2
$ws.Range("G69").Value2 = 2
$ws.Range("H68").Value2 = 4
$ws.Range("H69").Value2 = 1
$ws.Range("I68").Value2 = 2
$ws.Range("I69").Value2 = 0
$ws.Range("J68").Value2 = 2
$ws.Range("J69").Value2 = 1
$ws.Range("K68").Value2 = 'A'
$ws.Range("K69").Value2 = 'H'
$ws.Range("L68").Value2 = 4
$ws.Range("L69").Value2 = 1.909
$ws.Range("M68").Value2 = 4
$ws.Range("M69").Value2 = 3.8
$ws.Range("N68").Value2 = 1.65
$ws.Range("N69").Value2 = 3.25
$ws.Range("O68").Value2 = 4
$ws.Range("O69").Value2 = 1.909
$ws.Range("P68").Value2 = 3.75
$ws.Range("P69").Value2 = 3.75
$ws.Range("Q68").Value2 = 1.666
$ws.Range("Q69").Value2 = 3.4
$ws.Range("R68").Value2 = 0.75
$ws.Range("R69").Value2 = -0.5
$ws.Range("S68").Value2 = 1.925
$ws.Range("S69").Value2 = 1.925
$ws.Range("T68").Value2 = 1.875
$ws.Range("T69").Value2 = 1.875
$ws.Range("U68").Value2 = 3
$ws.Range("U69").Value2 = 3
$ws.Range("V68").Value2 = 1.825
$ws.Range("V69").Value2 = 1.825
$ws.Range("W68").Value2 = 1.975
$ws.Range("W69").Value2 = 1.975
$ws.Range("X68").Value2 = -1
$ws.Range("X69").Value2 = 0.909
$ws.Range("Y68").Value2 = -1
$ws.Range("Y69").Value2 = -1
$ws.Range("Z68").Value2 = 0.6659999999999999
$ws.Range("Z69").Value2 = -1
$ws.Range("AA68").Value2 = -1
$ws.Range("AA69").Value2 = 0.925
$ws.Range("AB68").Value2 = 0.875
$ws.Range("AB69").Value2 = -1
$ws.Range("AC68").Value2 = 0.825
$ws.Range("AC69").Value2 = 0
$ws.Range("AD68").Value2 = -1
$ws.Range("AD69").Value2 = 0

# Swap row 118 and row 119 (columns B:AD), keep column A (index) unchanged
$ws.Range("B118").Value2 = 7818294
$ws.Range("B119").Value2 = 7818846
$ws.Range("C118").Value2 = 'Qatar Stars League'
$ws.Range("C119").Value2 = 'Qatar Stars League'
$ws.Range("D118").Value2 = 45399.52083333334
$ws.Range("D119").Value2 = 45399.52083333334
$ws.Range("E118").Value2 = 'Umm Salal'
$ws.Range("E119").Value2 = 'AlMuaidar'
$ws.Range("F118").Value2 = 'AlArabi Doha'
$ws.Range("F119").Value2 = 'Al Duhail'
$ws.Range("G118").Value2 = 2
$ws.Range("G119").Value2 = 2
$ws.Range("H118").Value2 = 1
$ws.Range("H119").Value2 = 5
$ws.Range("I118").Value2 = 0
$ws.Range("I119").Value2 = 2
$ws.Range("J118").Value2 = 1
$ws.Range("J119").Value2 = 2
$ws.Range("K118").Value2 = 'H'
$ws.Range("K119").Value2 = 'A'
$ws.Range("L118").Value2 = 3.75
$ws.Range("L119").Value2 = 4
$ws.Range("M118").Value2 = 3.6
$ws.Range("M119").Value2 = 4
$ws.Range("N118").Value2 = 1.8
$ws.Range("N119").Value2 = 1.666
$ws.Range("O118").Value2 = 4.333
$ws.Range("O119").Value2 = 4
$ws.Range("P118").Value2 = 3.8
$ws.Range("P119").Value2 = 3.8
$ws.Range("Q118").Value2 = 1.65
$ws.Range("Q119").Value2 = 1.7
$ws.Range("R118").Value2 = 0.75
$ws.Range("R119").Value2 = 0.75
$ws.Range("S118").Value2 = 1.95
$ws.Range("S119").Value2 = 1.9
$ws.Range("T118").Value2 = 1.85
$ws.Range("T119").Value2 = 1.9
$ws.Range("U118").Value2 = 3.25
$ws.Range("U119").Value2 = 3.25
$ws.Range("V118").Value2 = 1.975
$ws.Range("V119").Value2 = 2
$ws.Range("W118").Value2 = 1.825
$ws.Range("W119").Value2 = 1.8
$ws.Range("X118").Value2 = 3.333
$ws.Range("X119").Value2 = -1
$ws.Range("Y118").Value2 = -1
$ws.Range("Y119").Value2 = -1
$ws.Range("Z118").Value2 = -1
$ws.Range("Z119").Value2 = 0.7
$ws.Range("AA118").Value2 = 0.95
$ws.Range("AA119").Value2 = -1
$ws.Range("AB118").Value2 = -1
$ws.Range("AB119").Value2 = 0.8999999999999999
$ws.Range("AC118").Value2 = -0.5
$ws.Range("AC119").Value2 = 1
$ws.Range("AD118").Value2 = 0.4125
$ws.Range("AD119").Value2 = -1

# Swap row 123 and row 125 (columns B:AD), keep column A (index) unchanged
$ws.Range("B123").Value2 = 7004663
$ws.Range("B125").Value2 = 7004661
$ws.Range("C123").Value2 = 'Qatar Stars League'
$ws.Range("C125").Value2 = 'Qatar Stars League'
$ws.Range("D123").Value2 = 45406.52083333334
$ws.Range("D125").Value2 = 45406.52083333334
$ws.Range("E123").Value2 = 'AlWakrah SC'
$ws.Range("E125").Value2 = 'Umm Salal'
$ws.Range("F123").Value2 = 'Al Gharafa'
$ws.Range("F125").Value2 = 'Al Duhail'
$ws.Range("G123").Value2 = 2
$ws.Range("G125").Value2 = 2
$ws.Range("H123").Value2 = 4
$ws.Range("H125").Value2 = 0
$ws.Range("I123").Value2 = 0
$ws.Range("I125").Value2 = 0
$ws.Range("J123").Value2 = 1
$ws.Range("J125").Value2 = 0
$ws.Range("K123").Value2 = 'A'
$ws.Range("K125").Value2 = 'H'
$ws.Range("L123").Value2 = 2.375
$ws.Range("L125").Value2 = 4.75
$ws.Range("M123").Value2 = 3.5
$ws.Range("M125").Value2 = 4.333
$ws.Range("N123").Value2 = 2.6
$ws.Range("N125").Value2 = 1.533
$ws.Range("O123").Value2 = 3.3
$ws.Range("O125").Value2 = 5.5
$ws.Range("P123").Value2 = 3.6
$ws.Range("P125").Value2 = 4.5
$ws.Range("Q123").Value2 = 1.95
$ws.Range("Q125").Value2 = 1.45
$ws.Range("R123").Value2 = 0.5
$ws.Range("R125").Value2 = 1.25
$ws.Range("S123").Value2 = 1.8
$ws.Range("S125").Value2 = 1.85
$ws.Range("T123").Value2 = 2
$ws.Range("T125").Value2 = 1.95
$ws.Range("U123").Value2 = 3.25
$ws.Range("U125").Value2 = 3.25
$ws.Range("V123").Value2 = 1.975
$ws.Range("V125").Value2 = 1.875
$ws.Range("W123").Value2 = 1.825
$ws.Range("W125").Value2 = 1.925
$ws.Range("X123").Value2 = -1
$ws.Range("X125").Value2 = 4.5
$ws.Range("Y123").Value2 = -1
$ws.Range("Y125").Value2 = -1
$ws.Range("Z123").Value2 = 0.95
$ws.Range("Z125").Value2 = -1
$ws.Range("AA123").Value2 = -1
$ws.Range("AA125").Value2 = 0.8500000000000001
$ws.Range("AB123").Value2 = 1
$ws.Range("AB125").Value2 = -1
$ws.Range("AC123").Value2 = 0.9750000000000001
$ws.Range("AC125").Value2 = -1
$ws.Range("AD123").Value2 = -1
$ws.Range("AD125").Value2 = 0.925
